{"js": "// Update the multiplication-answer table: each text cell holds a string\n// like \"92\u00d754=4968\" that must become a new equation (per the commit diff).\n// Cells are addressed by (row, col) rather than by searching for the old\n// text, because some old values repeat (e.g. \"31\u00d785=2635\" appears twice\n// but maps to two different replacements depending on position).\nconst replacements = [\n  { row: 0, col: 0, oldText: \"92\u00d754=4968\", newText: \"70\u00d768=4760\" },\n  { row: 0, col: 1, oldText: \"26\u00d754=1404\", newText: \"21\u00d743=903\" },\n  { row: 0, col: 2, oldText: \"50\u00d766=3300\", newText: \"36\u00d743=1548\" },\n  { row: 0, col: 3, oldText: \"47\u00d712=564\", newText: \"48\u00d761=2928\" },\n  { row: 0, col: 4, oldText: \"52\u00d760=3120\", newText: \"25\u00d746=1150\" },\n\n  { row: 4, col: 0, oldText: \"13\u00d774=962\", newText: \"88\u00d747=4136\" },\n  { row: 4, col: 1, oldText: \"95\u00d745=4275\", newText: \"32\u00d784=2688\" },\n  { row: 4, col: 2, oldText: \"82\u00d748=3936\", newText: \"96\u00d774=7104\" },\n  { row: 4, col: 3, oldText: \"69\u00d761=4209\", newText: \"27\u00d749=1323\" },\n  { row: 4, col: 4, oldText: \"43\u00d719=817\", newText: \"48\u00d740=1920\" },\n\n  { row: 9, col: 0, oldText: \"66\u00d749=3234\", newText: \"93\u00d788=8184\" },\n  { row: 9, col: 1, oldText: \"93\u00d759=5487\", newText: \"68\u00d760=4080\" },\n  { row: 9, col: 2, oldText: \"62\u00d726=1612\", newText: \"87\u00d762=5394\" },\n  { row: 9, col: 3, oldText: \"57\u00d725=1425\", newText: \"22\u00d743=946\" },\n  { row: 9, col: 4, oldText: \"31\u00d785=2635\", newText: \"74\u00d759=4366\" },\n\n  { row: 14, col: 0, oldText: \"99\u00d735=3465\", newText: \"82\u00d752=4264\" },\n  { row: 14, col: 1, oldText: \"31\u00d785=2635\", newText: \"88\u00d735=3080\" },\n  { row: 14, col: 2, oldText: \"90\u00d769=6210\", newText: \"94\u00d780=7520\" },\n  { row: 14, col: 3, oldText: \"86\u00d755=4730\", newText: \"39\u00d725=975\" },\n  { row: 14, col: 4, oldText: \"26\u00d775=1950\", newText: \"26\u00d795=2470\" },\n\n  { row: 19, col: 0, oldText: \"99\u00d765=6435\", newText: \"46\u00d752=2392\" },\n  { row: 19, col: 1, oldText: \"61\u00d766=4026\", newText: \"66\u00d716=1056\" },\n  { row: 19, col: 2, oldText: \"65\u00d737=2405\", newText: \"11\u00d714=154\" },\n  { row: 19, col: 3, oldText: \"93\u00d712=1116\", newText: \"35\u00d746=1610\" },\n  { row: 19, col: 4, oldText: \"74\u00d736=2664\", newText: \"79\u00d766=5214\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Load the current cell values first so we can sanity-check we are about to\n// overwrite the expected old equation before mutating anything.\nconst cells = replacements.map((r) => table.getCell(r.row, r.col));\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\nreplacements.forEach((r, i) => {\n  const cell = cells[i];\n  const current = (cell.value || \"\").trim();\n  if (current === r.oldText) {\n    cell.value = r.newText;\n  } else {\n    // Fallback: the cell didn't hold the text we expected (e.g. different\n    // whitespace/newline handling) \u2014 replace via the cell body's range so\n    // the edit still lands on the intended equation text.\n    cell.body.search(r.oldText, { matchCase: true }).items.forEach((rng) => {\n      rng.insertText(r.newText, \"Replace\");\n    });\n  }\n});\nawait context.sync();\n", "ps1": "# Update the multiplication-answer table: each text cell holds a string\n# like \"92x54=4968\" that must become a new equation (per the commit diff).\n# Cells are addressed by (row, col) -- 1-based, as in the Word COM object\n# model -- rather than by searching for the old text, because some old\n# values repeat (e.g. \"31x85=2635\" appears twice but maps to two different\n# replacements depending on position).\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1;  Col = 1; Old = \"92\u00d754=4968\"; New = \"70\u00d768=4760\" },\n    @{ Row = 1;  Col = 2; Old = \"26\u00d754=1404\"; New = \"21\u00d743=903\" },\n    @{ Row = 1;  Col = 3; Old = \"50\u00d766=3300\"; New = \"36\u00d743=1548\" },\n    @{ Row = 1;  Col = 4; Old = \"47\u00d712=564\"; New = \"48\u00d761=2928\" },\n    @{ Row = 1;  Col = 5; Old = \"52\u00d760=3120\"; New = \"25\u00d746=1150\" },\n\n    @{ Row = 5;  Col = 1; Old = \"13\u00d774=962\"; New = \"88\u00d747=4136\" },\n    @{ Row = 5;  Col = 2; Old = \"95\u00d745=4275\"; New = \"32\u00d784=2688\" },\n    @{ Row = 5;  Col = 3; Old = \"82\u00d748=3936\"; New = \"96\u00d774=7104\" },\n    @{ Row = 5;  Col = 4; Old = \"69\u00d761=4209\"; New = \"27\u00d749=1323\" },\n    @{ Row = 5;  Col = 5; Old = \"43\u00d719=817\"; New = \"48\u00d740=1920\" },\n\n    @{ Row = 10; Col = 1; Old = \"66\u00d749=3234\"; New = \"93\u00d788=8184\" },\n    @{ Row = 10; Col = 2; Old = \"93\u00d759=5487\"; New = \"68\u00d760=4080\" },\n    @{ Row = 10; Col = 3; Old = \"62\u00d726=1612\"; New = \"87\u00d762=5394\" },\n    @{ Row = 10; Col = 4; Old = \"57\u00d725=1425\"; New = \"22\u00d743=946\" },\n    @{ Row = 10; Col = 5; Old = \"31\u00d785=2635\"; New = \"74\u00d759=4366\" },\n\n    @{ Row = 15; Col = 1; Old = \"99\u00d735=3465\"; New = \"82\u00d752=4264\" },\n    @{ Row = 15; Col = 2; Old = \"31\u00d785=2635\"; New = \"88\u00d735=3080\" },\n    @{ Row = 15; Col = 3; Old = \"90\u00d769=6210\"; New = \"94\u00d780=7520\" },\n    @{ Row = 15; Col = 4; Old = \"86\u00d755=4730\"; New = \"39\u00d725=975\" },\n    @{ Row = 15; Col = 5; Old = \"26\u00d775=1950\"; New = \"26\u00d795=2470\" },\n\n    @{ Row = 20; Col = 1; Old = \"99\u00d765=6435\"; New = \"46\u00d752=2392\" },\n    @{ Row = 20; Col = 2; Old = \"61\u00d766=4026\"; New = \"66\u00d716=1056\" },\n    @{ Row = 20; Col = 3; Old = \"65\u00d737=2405\"; New = \"11\u00d714=154\" },\n    @{ Row = 20; Col = 4; Old = \"93\u00d712=1116\"; New = \"35\u00d746=1610\" },\n    @{ Row = 20; Col = 5; Old = \"74\u00d736=2664\"; New = \"79\u00d766=5214\" }\n)\n\nforeach ($r in $replacements) {\n    $cell = $table.Cell($r.Row, $r.Col)\n    # Cell.Range.Text includes the trailing paragraph/cell-end marks (CR +\n    # BEL); strip them before comparing against the expected old equation.\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $r.Old) {\n        # Fallback: locate the exact old text within the cell and replace\n        # just that text, in case the direct (row, col) read didn't line up\n        # with what we expected.\n        $findRange = $cell.Range\n        $findRange.Find.ClearFormatting()\n        $findRange.Find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n    } else {\n        $cell.Range.Text = $r.New\n    }\n}\n"}
